$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 437.81396
$ws.Range("I15").Value = 437.81396
$ws.Range("K15").Value = 1313.44188
$ws.Range("M15").Value = -1144.44188

$ws.Range("H46").Value = 4492.857
$ws.Range("I46").Value = 3650
$ws.Range("J46").Value = 4633.3335
$ws.Range("K46").Value = 10950
$ws.Range("L46").Value = 13900.0005
$ws.Range("M46").Value = -10831
$ws.Range("N46").Value = -14138.0005

$ws.Range("H60").Value = 4492.857
$ws.Range("I60").Value = 3650
$ws.Range("J60").Value = 4633.3335
$ws.Range("K60").Value = 10950
$ws.Range("L60").Value = 13900.0005
$ws.Range("M60").Value = -10466
$ws.Range("N60").Value = -14868.0005

$ws.Range("H74").Value = 4868.952
$ws.Range("I74").Value = 3660.889
$ws.Range("K74").Value = 3660.889
$ws.Range("M74").Value = -2724.889

$ws.Range("H77").Value = 4868.952
$ws.Range("I77").Value = 3660.889
$ws.Range("K77").Value = 18304.445
$ws.Range("M77").Value = -13624.445

$ws.Range("H132").Value = 8145.8
$ws.Range("I132").Value = 8698.929
$ws.Range("J132").Value = 402
$ws.Range("K132").Value = 26096.787
$ws.Range("L132").Value = 1206
$ws.Range("M132").Value = -23566.787
$ws.Range("N132").Value = -6266

$ws.Range("H137").Value = 3751.6943
$ws.Range("I137").Value = 2079.6538
$ws.Range("K137").Value = 6238.9614
$ws.Range("M137").Value = -3688.9614

$ws.Range("H138").Value = 4986.2256
$ws.Range("I138").Value = 1902.4
$ws.Range("K138").Value = 5707.200000000001
$ws.Range("M138").Value = -567.2000000000007

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 55556596
$ws.Range("I122").Value = 1258.25
$ws.Range("J122").Value = 166667260
$ws.Range("K122").Value = 3774.75
$ws.Range("L122").Value = 500001780
$ws.Range("M122").Value = -1324.75
$ws.Range("N122").Value = -500006680

$ws.Range("H135").Value = 105362.836
$ws.Range("J135").Value = 105362.836
$ws.Range("L135").Value = 105362.836
$ws.Range("N135").Value = -115502.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1915.0476
$ws.Range("I20").Value = 1293.625
$ws.Range("K20").Value = 1293.625
$ws.Range("M20").Value = -1046.625

$ws.Range("H75").Value = 30356.6
$ws.Range("I75").Value = 5445.75
$ws.Range("K75").Value = 5445.75
$ws.Range("M75").Value = -4509.75

$ws.Range("H78").Value = 30356.6
$ws.Range("I78").Value = 5445.75
$ws.Range("K78").Value = 16337.25
$ws.Range("M78").Value = -11657.25

$ws.Range("H86").Value = 32391.9
$ws.Range("I86").Value = 18402.166
$ws.Range("J86").Value = 53376.5
$ws.Range("K86").Value = 18402.166
$ws.Range("L86").Value = 53376.5
$ws.Range("M86").Value = -17279.166
$ws.Range("N86").Value = -55622.5

$ws.Range("H89").Value = 32391.9
$ws.Range("I89").Value = 18402.166
$ws.Range("J89").Value = 53376.5
$ws.Range("K89").Value = 92010.83
$ws.Range("L89").Value = 266882.5
$ws.Range("M89").Value = -86394.83
$ws.Range("N89").Value = -278114.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 8992.5
$ws.Range("I22").Value = 12862.25
$ws.Range("K22").Value = 12862.25
$ws.Range("M22").Value = -12512.25

$ws.Range("H58").Value = 1353.375
$ws.Range("I58").Value = 1419.8
$ws.Range("K58").Value = 1419.8
$ws.Range("M58").Value = -1216.8

$ws.Range("H133").Value = 84850.78999999999
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

$ws.Range("H136").Value = 1353.375
$ws.Range("I136").Value = 1419.8
$ws.Range("K136").Value = 4259.4
$ws.Range("M136").Value = -1709.4

$ws.Range("H141").Value = 96537.55
$ws.Range("J141").Value = 103291.3
$ws.Range("L141").Value = 103291.3
$ws.Range("N141").Value = -113651.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2832.3333
$ws.Range("I22").Value = 998.5
$ws.Range("J22").Value = 6500
$ws.Range("K22").Value = 2995.5
$ws.Range("L22").Value = 19500
$ws.Range("M22").Value = -2826.5
$ws.Range("N22").Value = -19838

$ws.Range("H27").Value = 2832.3333
$ws.Range("I27").Value = 998.5
$ws.Range("J27").Value = 6500
$ws.Range("K27").Value = 2995.5
$ws.Range("L27").Value = 19500
$ws.Range("M27").Value = -2893.5
$ws.Range("N27").Value = -19704

$ws.Range("H44").Value = 5795.8184
$ws.Range("I44").Value = 6766.8
$ws.Range("J44").Value = 2761.5
$ws.Range("K44").Value = 20300.4
$ws.Range("L44").Value = 8284.5
$ws.Range("M44").Value = -19902.4
$ws.Range("N44").Value = -9080.5

$ws.Range("H55").Value = 904904.8
$ws.Range("I55").Value = 2252924.5
$ws.Range("J55").Value = 6225
$ws.Range("K55").Value = 6758773.5
$ws.Range("L55").Value = 18675
$ws.Range("M55").Value = -6758596.5
$ws.Range("N55").Value = -19029

$ws.Range("H104").Value = 2995.2
$ws.Range("I104").Value = 2994
$ws.Range("J104").Value = 2996.4
$ws.Range("K104").Value = 8982
$ws.Range("L104").Value = 8989.200000000001
$ws.Range("M104").Value = -6361
$ws.Range("N104").Value = -14231.2

$ws.Range("H112").Value = 21555.445
$ws.Range("I112").Value = 24999
$ws.Range("J112").Value = 21125
$ws.Range("K112").Value = 74997
$ws.Range("L112").Value = 63375
$ws.Range("M112").Value = -73889
$ws.Range("N112").Value = -65591

$ws.Range("H123").Value = 1500
$ws.Range("I123").Value = 1500
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 4500
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -2050
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 6066178
$ws.Range("I132").Value = 2098.25
$ws.Range("J132").Value = 9531366
$ws.Range("K132").Value = 18884.25
$ws.Range("L132").Value = 85782294
$ws.Range("M132").Value = -16354.25
$ws.Range("N132").Value = -85787354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 967.75
$ws.Range("I107").Value = 656
$ws.Range("J107").Value = 1903
$ws.Range("K107").Value = 656
$ws.Range("L107").Value = 1903
$ws.Range("M107").Value = 1264
$ws.Range("N107").Value = -5743

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H132").Value = 6264
$ws.Range("I132").Value = 3396.5
$ws.Range("J132").Value = 11999
$ws.Range("K132").Value = 10189.5
$ws.Range("L132").Value = 35997
$ws.Range("M132").Value = -7659.5
$ws.Range("N132").Value = -41057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 39999.668
$ws.Range("J38").Value = 20000
$ws.Range("L38").Value = 20000
$ws.Range("N38").Value = -20820

$ws.Range("H132").Value = 3386.2222
$ws.Range("I132").Value = 3323.7144
$ws.Range("J132").Value = 3605
$ws.Range("K132").Value = 9971.143199999999
$ws.Range("L132").Value = 10815
$ws.Range("M132").Value = -7441.143199999999
$ws.Range("N132").Value = -15875

$ws.Range("H136").Value = 7099.933
$ws.Range("I136").Value = 2833.2222
$ws.Range("J136").Value = 13500
$ws.Range("K136").Value = 8499.6666
$ws.Range("L136").Value = 40500
$ws.Range("M136").Value = -5949.6666
$ws.Range("N136").Value = -45600
